$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Replace the four paragraphs:
#      "Die erstellten Projekt-WPS-Scripts sind hier abgelegt..."
#      "Öffentliche GitHub-URLs im Ablageordner auf dem BSCW ablegen! ..."
#      "Namenskonvention URL: M122_Klasse_Thema_Name_Name  "
#      "-> Ein Branch und separater Doku-Ordner pro Teammitglied erstellen"
#    with the single, reworded paragraph:
#      "Namenskonvention URL: M122_BI20a_GitGui_Afkhami_Greil"
# ---------------------------------------------------------------------------
$rngStart = $d.Content
$rngStart.Find.Execute("Die erstellten Projekt-WPS-Scripts", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start = $rngStart.Paragraphs(1).Range.Start

$rngEnd = $d.Content
$rngEnd.Find.Execute("pro Teammitglied erstellen", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$end = $rngEnd.Paragraphs(1).Range.End

$target1 = $d.Range($start, $end)

$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Textkrper"/><w:rPr><w:b/><w:i/><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="it-IT" w:eastAsia="de-DE"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="it-IT"/></w:rPr><w:t>Namenskonvention</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve"> URL: </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="4F81BD" w:themeColor="accent1"/><w:lang w:val="it-IT"/></w:rPr><w:t>M122_</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="4F81BD" w:themeColor="accent1"/><w:lang w:val="it-IT"/></w:rPr><w:t>BI20a</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="4F81BD" w:themeColor="accent1"/><w:lang w:val="it-IT"/></w:rPr><w:t>_</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="4F81BD" w:themeColor="accent1"/><w:lang w:val="it-IT"/></w:rPr><w:t>GitGui</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="4F81BD" w:themeColor="accent1"/><w:lang w:val="it-IT"/></w:rPr><w:t>_</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="4F81BD" w:themeColor="accent1"/><w:lang w:val="it-IT"/></w:rPr><w:t>Afkhami</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="4F81BD" w:themeColor="accent1"/><w:lang w:val="it-IT"/></w:rPr><w:t>_</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="4F81BD" w:themeColor="accent1"/><w:lang w:val="it-IT"/></w:rPr><w:t>G</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="4F81BD" w:themeColor="accent1"/><w:lang w:val="it-IT"/></w:rPr><w:t>reil</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# 2) Replace the three paragraphs:
#      "Aufgrund unten beschriebener Umstände sind Anpassungen ... worden:"
#      "..."
#      "Umstände / Anpassungen / Veränderungen"
#    with the single paragraph:
#      "GitGui.ps1:" followed by a line break
# ---------------------------------------------------------------------------
$rngStart2 = $d.Content
$rngStart2.Find.Execute("Aufgrund unten beschriebener", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start2 = $rngStart2.Paragraphs(1).Range.Start

$rngEnd2 = $d.Content
$rngEnd2.Find.Execute("Umstände / Anpassungen / Veränderungen", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$end2 = $rngEnd2.Paragraphs(1).Range.End

$target2 = $d.Range($start2, $end2)

$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:lang w:eastAsia="de-DE"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:eastAsia="de-DE"/></w:rPr><w:t>GitGui.ps1</w:t></w:r><w:r><w:rPr><w:lang w:eastAsia="de-DE"/></w:rPr><w:t>:</w:t></w:r><w:r><w:rPr><w:lang w:eastAsia="de-DE"/></w:rPr><w:br/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target2.InsertXML($xml2)
